# AWS bash cmd line
# Refresh the "last status check" timestamp banner and the MOL Olomoucka
# price row (B7/C7/D7/E7) with the latest scrape results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status-check banner in F1 (17:30 -> 17:45).
$ws.Range("F1").Value = "Last status check on: 21.01.2022 17:45"

# New scrape: price moved up from 36.7 to 36.9 (old price becomes 36.7).
$ws.Range("B7").Value = 36.9
$ws.Range("C7").Value = 36.7

# Delta is now written out as a plain text string with an explicit sign,
# rather than a numeric value.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "+0.2"
$ws.Range("D7").Style = "Normal"

# Timestamp of this particular row's check, written as literal text
# instead of a serial date.
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2022-01-21 17:45:21"
$ws.Range("E7").Style = "Normal"
